$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(2, 46056.01041666666, 104),
  @(3, 46056.02083333334, 98),
  @(4, 46056.03125, 94),
  @(5, 46056.04166666666, 91),
  @(6, 46056.05208333334, 85),
  @(7, 46056.0625, 82),
  @(8, 46056.07291666666, 76),
  @(9, 46056.08333333334, 68),
  @(10, 46056.09375, 63),
  @(11, 46056.10416666666, 58),
  @(12, 46056.11458333334, 53),
  @(13, 46056.125, 47),
  @(14, 46056.13541666666, 43),
  @(15, 46056.14583333334, 45),
  @(16, 46056.15625, 47),
  @(17, 46056.16666666666, 49),
  @(18, 46056.17708333334, 39),
  @(19, 46056.1875, 48),
  @(20, 46056.19791666666, 50),
  @(21, 46056.20833333334, 52),
  @(22, 46056.21875, 47),
  @(23, 46056.22916666666, 43),
  @(24, 46056.23958333334, 41),
  @(25, 46056.25, 40),
  @(26, 46056.26041666666, 39),
  @(27, 46056.27083333334, 42),
  @(28, 46056.28125, 43),
  @(29, 46056.29166666666, 44),
  @(30, 46056.30208333334, 46),
  @(31, 46056.3125, 45),
  @(32, 46056.32291666666, 42),
  @(33, 46056.33333333334, 40),
  @(34, 46056.34375, 43),
  @(35, 46056.35416666666, 45),
  @(36, 46056.36458333334, 44),
  @(37, 46056.375, 45),
  @(38, 46056.38541666666, 0),
  @(39, 46056.39583333334, 44),
  @(40, 46056.40625, 0),
  @(41, 46056.41666666666, 43),
  @(42, 46056.42708333334, 46),
  @(43, 46056.4375, 48),
  @(44, 46056.44791666666, 53),
  @(45, 46056.45833333334, 56),
  @(46, 46056.46875, 58),
  @(47, 46056.47916666666, 67),
  @(48, 46056.48958333334, 74),
  @(49, 46056.5, 79),
  @(50, 46056.51041666666, 0),
  @(51, 46056.52083333334, 0),
  @(52, 46056.53125, 0),
  @(53, 46056.54166666666, 0),
  @(54, 46056.55208333334, 0),
  @(55, 46056.5625, 0),
  @(56, 46056.57291666666, 0),
  @(57, 46056.58333333334, 0),
  @(58, 46056.59375, 0),
  @(59, 46056.60416666666, 0),
  @(60, 46056.61458333334, 0),
  @(61, 46056.625, 0),
  @(62, 46056.63541666666, 0),
  @(63, 46056.64583333334, 0),
  @(64, 46056.65625, 0),
  @(65, 46056.66666666666, 0),
  @(66, 46056.67708333334, 0),
  @(67, 46056.6875, 0),
  @(68, 46056.69791666666, 0),
  @(69, 46056.70833333334, 0),
  @(70, 46056.71875, 0),
  @(71, 46056.72916666666, 0),
  @(72, 46056.73958333334, 0),
  @(73, 46056.75, 0),
  @(74, 46056.76041666666, 0),
  @(75, 46056.77083333334, 0),
  @(76, 46056.78125, 0),
  @(77, 46056.79166666666, 0),
  @(78, 46056.80208333334, 0),
  @(79, 46056.8125, 0),
  @(80, 46056.82291666666, 0),
  @(81, 46056.83333333334, 0),
  @(82, 46056.84375, 0),
  @(83, 46056.85416666666, 0),
  @(84, 46056.86458333334, 0),
  @(85, 46056.875, 0),
  @(86, 46056.88541666666, 0),
  @(87, 46056.89583333334, 0),
  @(88, 46056.90625, 0),
  @(89, 46056.91666666666, 0),
  @(90, 46056.92708333334, 0),
  @(91, 46056.9375, 0),
  @(92, 46056.94791666666, 0),
  @(93, 46056.95833333334, 0),
  @(94, 46056.96875, 0),
  @(95, 46056.97916666666, 0),
  @(96, 46056.98958333334, 0),
  @(97, 46057, 0)
)

foreach ($item in $data) {
  $r = $item[0]
  $ws.Cells.Item($r, 1).Value = $item[1]
  $ws.Cells.Item($r, 2).Value = $item[2]
}

Write-Host "Done. A2=" $ws.Cells.Item(2,1).Value " B2=" $ws.Cells.Item(2,2).Value
Write-Host "A97=" $ws.Cells.Item(97,1).Value " B97=" $ws.Cells.Item(97,2).Value
